# Update cryptocurrency price/volume data per latest scrape (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.466.25"
$ws.Range("E2").Value = "  +0.53%  "

$ws.Range("D3").Value = "1.571.75"
$ws.Range("E3").Value = "  +0.37%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.53%  "

$ws.Range("E5").Value = "  -0.44%  "

$ws.Range("D6").Value = "'290.52"
$ws.Range("E6").Value = "  +0.52%  "

$ws.Range("D7").Value = "'0.3696"
$ws.Range("E7").Value = "  -0.98%  "

$ws.Range("D8").Value = "'49.94"
$ws.Range("E8").Value = "  +1.63%  "

$ws.Range("D9").Value = "'0.3397"
$ws.Range("E9").Value = "  +1.10%  "

$ws.Range("D10").Value = "'1.147"
$ws.Range("E10").Value = "  +2.99%  "

$ws.Range("D11").Value = "'0.07552"
$ws.Range("E11").Value = "  +1.88%  "

$ws.Range("E12").Value = "  -0.69%  "

$ws.Range("D13").Value = "'21.19"
$ws.Range("E13").Value = "  +2.45%  "

$ws.Range("D14").Value = "'6.018"
$ws.Range("E14").Value = "  +2.98%  "

$ws.Range("D15").Value = "'6.998"
$ws.Range("E15").Value = "  +2.38%  "

$ws.Range("D16").Value = "1.573.67"
$ws.Range("E16").Value = "  +0.04%  "

$ws.Range("D17").Value = "'0.00001123"
$ws.Range("E17").Value = "  +1.51%  "

$ws.Range("D18").Value = "'90.41"
$ws.Range("E18").Value = "  +1.71%  "

$ws.Range("E19").Value = "  +1.60%  "

$ws.Range("E20").Value = "  -0.53%  "

$ws.Range("D21").Value = "'6.364"
$ws.Range("E21").Value = "  +4.06%  "

$ws.Range("D22").Value = "'16.39"
$ws.Range("E22").Value = "  +1.46%  "

$ws.Range("E23").Value = "  +3.41%  "

$ws.Range("D24").Value = "22.491.10"
$ws.Range("E24").Value = "  +0.63%  "

$ws.Range("D25").Value = "'2.364"
$ws.Range("E25").Value = "  -0.19%  "

$ws.Range("D26").Value = "'2.648"
$ws.Range("E26").Value = "  +6.24%  "

$ws.Range("D27").Value = "'20.02"
$ws.Range("E27").Value = "  +1.08%  "

$ws.Range("D28").Value = "'149.79"
$ws.Range("E28").Value = "  +1.75%  "

$ws.Range("D29").Value = "'5.065"

$ws.Range("D30").Value = "'124.77"
$ws.Range("E30").Value = "  +0.47%  "

$ws.Range("D31").Value = "1.750.49"
$ws.Range("E31").Value = "  +0.67%  "

$ws.Range("D32").Value = "'1.064"
$ws.Range("E32").Value = "  +9.66%  "

$ws.Range("D33").Value = "'6.224"
$ws.Range("E33").Value = "  +6.13%  "

$ws.Range("D34").Value = "'2.016"
$ws.Range("E34").Value = "  +1.50%  "

$ws.Range("D35").Value = "'9.822"
$ws.Range("E35").Value = "  +2.08%  "

$ws.Range("D36").Value = "'0.08383"
$ws.Range("E36").Value = "  -0.05%  "

$ws.Range("D37").Value = "'0.02482"
$ws.Range("E37").Value = "  +1.50%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'1.349"
$ws.Range("E38").Value = "  -2.76%  "

$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "'0.2300"
$ws.Range("E39").Value = "  +2.49%  "

$ws.Range("D40").Value = "'0.06522"
$ws.Range("E40").Value = "  +3.00%  "

$ws.Range("D41").Value = "'5.433"

$ws.Range("D42").Value = "'11.27"
$ws.Range("E42").Value = "  +3.31%  "

$ws.Range("D43").Value = "'0.6224"
$ws.Range("E43").Value = "  +1.26%  "

$ws.Range("D44").Value = "'14.14"
$ws.Range("E44").Value = "  +2.68%  "

$ws.Range("D46").Value = "'3.794"

$ws.Range("D47").Value = "'0.5872"
$ws.Range("E47").Value = "  +2.62%  "

$ws.Range("D48").Value = "'2.066"
$ws.Range("E48").Value = "  +1.99%  "

$ws.Range("D49").Value = "'126.89"
$ws.Range("E49").Value = "  +1.31%  "

$ws.Range("D50").Value = "'1.234"
$ws.Range("E50").Value = "  +0.83%  "

$ws.Range("D51").Value = "'0.07300"
$ws.Range("E51").Value = "  +0.19%  "
